# Scenarios.xlsx - "Added a new JMeter scripts and list of scenarios"
#
# Sheet1 row 7 (the "Insert Place / Place List" scenario) grows two extra
# sub-bullets ("Insert Country" / "Insert State"), and a new blank-looking
# row 8 is appended below it that just carries a lone line-break string in
# column B (with a left border on A8) - this is what the diff's new shared
# string / new row / new styles all encode.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the brand-new row 8 first so its lone "\n" shared string is
# registered in the shared-string table *before* the edited B7 string -
# that reproduces the table order used upstream (new <si> inserted ahead
# of the modified "Insert Place..." entry).
$ws.Range("B8").Value = "`n"

# Extend the existing B7 text with two extra lines.
$ws.Range("B7").Value = "Insert Place`nInsert Country`nInsert State`nPlace List"

# Row heights grow to fit the extra wrapped lines (row7: 4 lines, row8: 2 lines).
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 28.8

# A8 stays empty but picks up a thin left border; B8 keeps the word-wrap
# formatting used by the rest of column B.
$ws.Range("A8").Borders.Item(7).LineStyle = 1
$ws.Range("B8").WrapText = $true

# Match the saved selection / scroll position from the edit.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select()
